$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 792.34  # H17: 665.01 -> 792.34
$ws.Cells.Item(17, 10).Value = 790.69794  # J17: 658.0625 -> 790.69794
$ws.Cells.Item(17, 12).Value = 2372.09382  # L17: 1974.1875 -> 2372.09382
$ws.Cells.Item(17, 14).Value = -2708.09382  # N17: -2310.1875 -> -2708.09382
$ws.Cells.Item(94, 8).Value = 44000  # H94: 13904 -> 44000
$ws.Cells.Item(94, 9).Value = 0  # I94: 4880 -> 0
$ws.Cells.Item(94, 10).Value = 44000  # J94: 50000 -> 44000
$ws.Cells.Item(94, 11).Value = 0  # K94: 4880 -> 0
$ws.Cells.Item(94, 12).Value = 44000  # L94: 50000 -> 44000
$ws.Cells.Item(94, 13).ClearContents()  # M94: -4429 -> (removed)
$ws.Cells.Item(94, 14).Value = -44902  # N94: -50902 -> -44902
$ws.Cells.Item(100, 8).Value = 6293.95  # H100: 6548.8887 -> 6293.95
$ws.Cells.Item(100, 10).Value = 9802.1  # J100: 11252.75 -> 9802.1
$ws.Cells.Item(100, 12).Value = 9802.1  # L100: 11252.75 -> 9802.1
$ws.Cells.Item(100, 14).Value = -10884.1  # N100: -12334.75 -> -10884.1
$ws.Cells.Item(103, 8).Value = 540.8333  # H103: 541.6667 -> 540.8333
$ws.Cells.Item(103, 9).Value = 540.8333  # I103: 560 -> 540.8333
$ws.Cells.Item(103, 10).Value = 0  # J103: 450 -> 0
$ws.Cells.Item(103, 11).Value = 1622.4999  # K103: 1680 -> 1622.4999
$ws.Cells.Item(103, 12).Value = 0  # L103: 1350 -> 0
$ws.Cells.Item(103, 13).Value = -1036.4999  # M103: -1094 -> -1036.4999
$ws.Cells.Item(103, 14).ClearContents()  # N103: -2522 -> (removed)
$ws.Cells.Item(116, 8).Value = 231638  # H116: 420171.4 -> 231638
$ws.Cells.Item(116, 9).Value = 7374.8335  # I116: 116451 -> 7374.8335
$ws.Cells.Item(116, 10).Value = 366195.9  # J116: 723891.8 -> 366195.9
$ws.Cells.Item(116, 11).Value = 7374.8335  # K116: 116451 -> 7374.8335
$ws.Cells.Item(116, 12).Value = 366195.9  # L116: 723891.8 -> 366195.9
$ws.Cells.Item(116, 13).Value = -3932.8335  # M116: -113009 -> -3932.8335
$ws.Cells.Item(116, 14).Value = -373079.9  # N116: -730775.8 -> -373079.9
$ws.Cells.Item(132, 8).Value = 65668.59  # H132: 69396 -> 65668.59
$ws.Cells.Item(132, 9).Value = 68278.57000000001  # I132: 74636.44 -> 68278.57000000001
$ws.Cells.Item(132, 10).Value = 19994  # J132: 13498 -> 19994
$ws.Cells.Item(132, 11).Value = 204835.71  # K132: 223909.32 -> 204835.71
$ws.Cells.Item(132, 12).Value = 59982  # L132: 40494 -> 59982
$ws.Cells.Item(132, 13).Value = -202305.71  # M132: -221379.32 -> -202305.71
$ws.Cells.Item(132, 14).Value = -65042  # N132: -45554 -> -65042
$ws.Cells.Item(138, 8).Value = 7786.171  # H138: 3321.158 -> 7786.171
$ws.Cells.Item(138, 9).Value = 12091.952  # I138: 0 -> 12091.952
$ws.Cells.Item(138, 10).Value = 3265.1  # J138: 3321.158 -> 3265.1
$ws.Cells.Item(138, 11).Value = 36275.856  # K138: 0 -> 36275.856
$ws.Cells.Item(138, 12).Value = 9795.299999999999  # L138: 9963.474 -> 9795.299999999999
$ws.Cells.Item(138, 13).Value = -31135.856  # M138: None -> -31135.856
$ws.Cells.Item(138, 14).Value = -20075.3  # N138: -20243.474 -> -20075.3
$ws.Cells.Item(141, 8).Value = 630.3333  # H141: 637.8 -> 630.3333
$ws.Cells.Item(141, 9).Value = 630.3333  # I141: 637.8 -> 630.3333
$ws.Cells.Item(141, 11).Value = 1890.9999  # K141: 1913.4 -> 1890.9999
$ws.Cells.Item(141, 13).Value = 3289.0001  # M141: 3266.6 -> 3289.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8065531  # H32: 8334368.5 -> 8065531
$ws.Cells.Item(32, 9).Value = 8475531  # I32: 8621660 -> 8475531
$ws.Cells.Item(32, 10).Value = 2200  # J32: 2900 -> 2200
$ws.Cells.Item(32, 11).Value = 8475531  # K32: 8621660 -> 8475531
$ws.Cells.Item(32, 12).Value = 2200  # L32: 2900 -> 2200
$ws.Cells.Item(32, 13).Value = -8475244  # M32: -8621373 -> -8475244
$ws.Cells.Item(32, 14).Value = -2774  # N32: -3474 -> -2774
$ws.Cells.Item(61, 8).Value = 1589766.2  # H61: 1963176.6 -> 1589766.2
$ws.Cells.Item(61, 9).Value = 1853560.6  # I61: 2224000.2 -> 1853560.6
$ws.Cells.Item(61, 10).Value = 6999.6665  # J61: 7000 -> 6999.6665
$ws.Cells.Item(61, 11).Value = 1853560.6  # K61: 2224000.2 -> 1853560.6
$ws.Cells.Item(61, 12).Value = 6999.6665  # L61: 7000 -> 6999.6665
$ws.Cells.Item(61, 13).Value = -1853348.6  # M61: -2223788.2 -> -1853348.6
$ws.Cells.Item(61, 14).Value = -7423.6665  # N61: -7424 -> -7423.6665
$ws.Cells.Item(97, 8).Value = 678.04346  # H97: 758.8261 -> 678.04346
$ws.Cells.Item(97, 9).Value = 678.04346  # I97: 702.4545000000001 -> 678.04346
$ws.Cells.Item(97, 10).Value = 0  # J97: 1999 -> 0
$ws.Cells.Item(97, 11).Value = 678.04346  # K97: 702.4545000000001 -> 678.04346
$ws.Cells.Item(97, 12).Value = 0  # L97: 1999 -> 0
$ws.Cells.Item(97, 13).Value = -182.04346  # M97: -206.4545000000001 -> -182.04346
$ws.Cells.Item(97, 14).ClearContents()  # N97: -2991 -> (removed)
$ws.Cells.Item(102, 8).Value = 28710.154  # H102: 24976.133 -> 28710.154
$ws.Cells.Item(102, 9).Value = 28710.154  # I102: 24976.133 -> 28710.154
$ws.Cells.Item(102, 11).Value = 28710.154  # K102: 24976.133 -> 28710.154
$ws.Cells.Item(102, 13).Value = -27088.154  # M102: -23354.133 -> -27088.154
$ws.Cells.Item(122, 8).Value = 2780.375  # H122: 2892.4 -> 2780.375
$ws.Cells.Item(122, 9).Value = 1433.1111  # I122: 1474.75 -> 1433.1111
$ws.Cells.Item(122, 11).Value = 4299.3333  # K122: 4424.25 -> 4299.3333
$ws.Cells.Item(122, 13).Value = -1849.3333  # M122: -1974.25 -> -1849.3333
$ws.Cells.Item(136, 8).Value = 1589766.2  # H136: 1963176.6 -> 1589766.2
$ws.Cells.Item(136, 9).Value = 1853560.6  # I136: 2224000.2 -> 1853560.6
$ws.Cells.Item(136, 10).Value = 6999.6665  # J136: 7000 -> 6999.6665
$ws.Cells.Item(136, 11).Value = 5560681.800000001  # K136: 6672000.600000001 -> 5560681.800000001
$ws.Cells.Item(136, 12).Value = 20998.9995  # L136: 21000 -> 20998.9995
$ws.Cells.Item(136, 13).Value = -5558131.800000001  # M136: -6669450.600000001 -> -5558131.800000001
$ws.Cells.Item(136, 14).Value = -26098.9995  # N136: -26100 -> -26098.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1270.5264  # H20: 1269 -> 1270.5264
$ws.Cells.Item(20, 10).Value = 1745.875  # J20: 1809.8572 -> 1745.875
$ws.Cells.Item(20, 12).Value = 1745.875  # L20: 1809.8572 -> 1745.875
$ws.Cells.Item(20, 14).Value = -2239.875  # N20: -2303.8572 -> -2239.875
$ws.Cells.Item(64, 8).Value = 375.66666  # H64: 442.25 -> 375.66666
$ws.Cells.Item(64, 9).Value = 245.5  # I64: 248 -> 245.5
$ws.Cells.Item(64, 10).Value = 636  # J64: 636.5 -> 636
$ws.Cells.Item(64, 11).Value = 245.5  # K64: 248 -> 245.5
$ws.Cells.Item(64, 12).Value = 636  # L64: 636.5 -> 636
$ws.Cells.Item(64, 13).Value = -20.5  # M64: -23 -> -20.5
$ws.Cells.Item(64, 14).Value = -1086  # N64: -1086.5 -> -1086
$ws.Cells.Item(67, 8).Value = 375.66666  # H67: 442.25 -> 375.66666
$ws.Cells.Item(67, 9).Value = 245.5  # I67: 248 -> 245.5
$ws.Cells.Item(67, 10).Value = 636  # J67: 636.5 -> 636
$ws.Cells.Item(67, 11).Value = 245.5  # K67: 248 -> 245.5
$ws.Cells.Item(67, 12).Value = 636  # L67: 636.5 -> 636
$ws.Cells.Item(67, 13).Value = 534.5  # M67: 532 -> 534.5
$ws.Cells.Item(67, 14).Value = -2196  # N67: -2196.5 -> -2196
$ws.Cells.Item(86, 8).Value = 2080.611  # H86: 2153.4375 -> 2080.611
$ws.Cells.Item(86, 9).Value = 1904.8182  # I86: 1995.2222 -> 1904.8182
$ws.Cells.Item(86, 11).Value = 1904.8182  # K86: 1995.2222 -> 1904.8182
$ws.Cells.Item(86, 13).Value = -781.8181999999999  # M86: -872.2221999999999 -> -781.8181999999999
$ws.Cells.Item(89, 8).Value = 2080.611  # H89: 2153.4375 -> 2080.611
$ws.Cells.Item(89, 9).Value = 1904.8182  # I89: 1995.2222 -> 1904.8182
$ws.Cells.Item(89, 11).Value = 9524.091  # K89: 9976.110999999999 -> 9524.091
$ws.Cells.Item(89, 13).Value = -3908.091  # M89: -4360.110999999999 -> -3908.091
$ws.Cells.Item(105, 8).Value = 2281.8572  # H105: 2462.1667 -> 2281.8572
$ws.Cells.Item(105, 9).Value = 1649.75  # I105: 1799.6666 -> 1649.75
$ws.Cells.Item(105, 11).Value = 1649.75  # K105: 1799.6666 -> 1649.75
$ws.Cells.Item(105, 13).Value = 97.25  # M105: -52.66660000000002 -> 97.25
$ws.Cells.Item(134, 8).Value = 424940.7  # H134: 431253.22 -> 424940.7
$ws.Cells.Item(134, 10).Value = 275172.4  # J134: 292245.7 -> 275172.4
$ws.Cells.Item(134, 12).Value = 825517.2000000001  # L134: 876737.1000000001 -> 825517.2000000001
$ws.Cells.Item(134, 14).Value = -830587.2000000001  # N134: -881807.1000000001 -> -830587.2000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 85622  # H31: 97057.92 -> 85622
$ws.Cells.Item(31, 9).Value = 143565.4  # I31: 166074.8 -> 143565.4
$ws.Cells.Item(31, 10).Value = 21884.25  # J31: 24206.777 -> 21884.25
$ws.Cells.Item(31, 11).Value = 143565.4  # K31: 166074.8 -> 143565.4
$ws.Cells.Item(31, 12).Value = 21884.25  # L31: 24206.777 -> 21884.25
$ws.Cells.Item(31, 13).Value = -143270.4  # M31: -165779.8 -> -143270.4
$ws.Cells.Item(31, 14).Value = -22474.25  # N31: -24796.777 -> -22474.25
$ws.Cells.Item(34, 8).Value = 85622  # H34: 97057.92 -> 85622
$ws.Cells.Item(34, 9).Value = 143565.4  # I34: 166074.8 -> 143565.4
$ws.Cells.Item(34, 10).Value = 21884.25  # J34: 24206.777 -> 21884.25
$ws.Cells.Item(34, 11).Value = 143565.4  # K34: 166074.8 -> 143565.4
$ws.Cells.Item(34, 12).Value = 21884.25  # L34: 24206.777 -> 21884.25
$ws.Cells.Item(34, 13).Value = -143363.4  # M34: -165872.8 -> -143363.4
$ws.Cells.Item(34, 14).Value = -22288.25  # N34: -24610.777 -> -22288.25
$ws.Cells.Item(58, 8).Value = 460180.66  # H58: 443756.34 -> 460180.66
$ws.Cells.Item(58, 9).Value = 618768.9399999999  # I58: 589351.4 -> 618768.9399999999
$ws.Cells.Item(58, 10).Value = 7071.2856  # J58: 6971.2856 -> 7071.2856
$ws.Cells.Item(58, 11).Value = 618768.9399999999  # K58: 589351.4 -> 618768.9399999999
$ws.Cells.Item(58, 12).Value = 7071.2856  # L58: 6971.2856 -> 7071.2856
$ws.Cells.Item(58, 13).Value = -618565.9399999999  # M58: -589148.4 -> -618565.9399999999
$ws.Cells.Item(58, 14).Value = -7477.2856  # N58: -7377.2856 -> -7477.2856
$ws.Cells.Item(136, 8).Value = 460180.66  # H136: 443756.34 -> 460180.66
$ws.Cells.Item(136, 9).Value = 618768.9399999999  # I136: 589351.4 -> 618768.9399999999
$ws.Cells.Item(136, 10).Value = 7071.2856  # J136: 6971.2856 -> 7071.2856
$ws.Cells.Item(136, 11).Value = 1856306.82  # K136: 1768054.2 -> 1856306.82
$ws.Cells.Item(136, 12).Value = 21213.8568  # L136: 20913.8568 -> 21213.8568
$ws.Cells.Item(136, 13).Value = -1853756.82  # M136: -1765504.2 -> -1853756.82
$ws.Cells.Item(136, 14).Value = -26313.8568  # N136: -26013.8568 -> -26313.8568

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(62, 8).Value = 6995  # H62: 0 -> 6995
$ws.Cells.Item(62, 10).Value = 6995  # J62: 0 -> 6995
$ws.Cells.Item(62, 12).Value = 20985  # L62: 0 -> 20985
$ws.Cells.Item(62, 14).Value = -22357  # N62: None -> -22357
$ws.Cells.Item(65, 8).Value = 6995  # H65: 0 -> 6995
$ws.Cells.Item(65, 10).Value = 6995  # J65: 0 -> 6995
$ws.Cells.Item(65, 12).Value = 62955  # L65: 0 -> 62955
$ws.Cells.Item(65, 14).Value = -69819  # N65: None -> -69819
$ws.Cells.Item(113, 8).Value = 3183  # H113: 1560.6086 -> 3183
$ws.Cells.Item(113, 10).Value = 3183  # J113: 1560.6086 -> 3183
$ws.Cells.Item(113, 12).Value = 9549  # L113: 4681.825800000001 -> 9549
$ws.Cells.Item(113, 14).Value = -13889  # N113: -9021.825800000001 -> -13889

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 11369669  # H132: 11241961 -> 11369669
$ws.Cells.Item(132, 9).Value = 15879172  # I132: 15631117 -> 15879172
$ws.Cells.Item(132, 10).Value = 5721.6  # J132: 5721.88 -> 5721.6
$ws.Cells.Item(132, 11).Value = 47637516  # K132: 46893351 -> 47637516
$ws.Cells.Item(132, 12).Value = 17164.8  # L132: 17165.64 -> 17164.8
$ws.Cells.Item(132, 13).Value = -47634986  # M132: -46890821 -> -47634986
$ws.Cells.Item(132, 14).Value = -22224.8  # N132: -22225.64 -> -22224.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3497.7727  # H7: 3482.1428 -> 3497.7727
$ws.Cells.Item(7, 9).Value = 3386.1667  # I7: 3395.8333 -> 3386.1667
$ws.Cells.Item(7, 11).Value = 3386.1667  # K7: 3395.8333 -> 3386.1667
$ws.Cells.Item(7, 13).Value = -3274.1667  # M7: -3283.8333 -> -3274.1667
$ws.Cells.Item(13, 8).Value = 1500  # H13: 20000 -> 1500
$ws.Cells.Item(13, 9).Value = 1500  # I13: 20000 -> 1500
$ws.Cells.Item(13, 11).Value = 1500  # K13: 20000 -> 1500
$ws.Cells.Item(13, 13).Value = -1360  # M13: -19860 -> -1360
$ws.Cells.Item(22, 8).Value = 897.7143  # H22: 859.13336 -> 897.7143
$ws.Cells.Item(22, 9).Value = 387.8  # I22: 381.54544 -> 387.8
$ws.Cells.Item(22, 11).Value = 387.8  # K22: 381.54544 -> 387.8
$ws.Cells.Item(22, 13).Value = -92.80000000000001  # M22: -86.54543999999999 -> -92.80000000000001
$ws.Cells.Item(27, 8).Value = 897.7143  # H27: 859.13336 -> 897.7143
$ws.Cells.Item(27, 9).Value = 387.8  # I27: 381.54544 -> 387.8
$ws.Cells.Item(27, 11).Value = 387.8  # K27: 381.54544 -> 387.8
$ws.Cells.Item(27, 13).Value = -280.8  # M27: -274.54544 -> -280.8
$ws.Cells.Item(74, 8).Value = 69150.125  # H74: 70426.42999999999 -> 69150.125
$ws.Cells.Item(74, 10).Value = 71900.14  # J74: 73847.5 -> 71900.14
$ws.Cells.Item(74, 12).Value = 71900.14  # L74: 73847.5 -> 71900.14
$ws.Cells.Item(74, 14).Value = -73896.14  # N74: -75843.5 -> -73896.14
$ws.Cells.Item(77, 8).Value = 69150.125  # H77: 70426.42999999999 -> 69150.125
$ws.Cells.Item(77, 10).Value = 71900.14  # J77: 73847.5 -> 71900.14
$ws.Cells.Item(77, 12).Value = 215700.42  # L77: 221542.5 -> 215700.42
$ws.Cells.Item(77, 14).Value = -225684.42  # N77: -231526.5 -> -225684.42
$ws.Cells.Item(126, 8).Value = 3497.7727  # H126: 3482.1428 -> 3497.7727
$ws.Cells.Item(126, 9).Value = 3386.1667  # I126: 3395.8333 -> 3386.1667
$ws.Cells.Item(126, 11).Value = 10158.5001  # K126: 10187.4999 -> 10158.5001
$ws.Cells.Item(126, 13).Value = -7688.500100000001  # M126: -7717.499899999999 -> -7688.500100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 555  # H18: 1018.3333 -> 555
$ws.Cells.Item(18, 10).Value = 555  # J18: 1018.3333 -> 555
$ws.Cells.Item(18, 12).Value = 555  # L18: 1018.3333 -> 555
$ws.Cells.Item(18, 14).Value = -901  # N18: -1364.3333 -> -901
$ws.Cells.Item(33, 8).Value = 6999  # H33: 0 -> 6999
$ws.Cells.Item(33, 9).Value = 6988  # I33: 0 -> 6988
$ws.Cells.Item(33, 10).Value = 7021  # J33: 0 -> 7021
$ws.Cells.Item(33, 11).Value = 6988  # K33: 0 -> 6988
$ws.Cells.Item(33, 12).Value = 7021  # L33: 0 -> 7021
$ws.Cells.Item(33, 13).Value = -6738  # M33: None -> -6738
$ws.Cells.Item(33, 14).Value = -7521  # N33: None -> -7521
$ws.Cells.Item(36, 8).Value = 6999  # H36: 0 -> 6999
$ws.Cells.Item(36, 9).Value = 6988  # I36: 0 -> 6988
$ws.Cells.Item(36, 10).Value = 7021  # J36: 0 -> 7021
$ws.Cells.Item(36, 11).Value = 6988  # K36: 0 -> 6988
$ws.Cells.Item(36, 12).Value = 7021  # L36: 0 -> 7021
$ws.Cells.Item(36, 13).Value = -6738  # M36: None -> -6738
$ws.Cells.Item(36, 14).Value = -7521  # N36: None -> -7521
$ws.Cells.Item(81, 8).Value = 6150.3125  # H81: 6360.4 -> 6150.3125
$ws.Cells.Item(81, 10).Value = 6400.778  # J81: 6826 -> 6400.778
$ws.Cells.Item(81, 12).Value = 12801.556  # L81: 13652 -> 12801.556
$ws.Cells.Item(81, 14).Value = -14923.556  # N81: -15774 -> -14923.556
$ws.Cells.Item(84, 8).Value = 6150.3125  # H84: 6360.4 -> 6150.3125
$ws.Cells.Item(84, 10).Value = 6400.778  # J84: 6826 -> 6400.778
$ws.Cells.Item(84, 12).Value = 64007.78  # L84: 68260 -> 64007.78
$ws.Cells.Item(84, 14).Value = -74615.78  # N84: -78868 -> -74615.78
$ws.Cells.Item(87, 8).Value = 99999  # H87: 0 -> 99999
$ws.Cells.Item(87, 10).Value = 99999  # J87: 0 -> 99999
$ws.Cells.Item(87, 12).Value = 99999  # L87: 0 -> 99999
$ws.Cells.Item(87, 14).Value = -102495  # N87: None -> -102495
$ws.Cells.Item(90, 8).Value = 99999  # H90: 0 -> 99999
$ws.Cells.Item(90, 10).Value = 99999  # J90: 0 -> 99999
$ws.Cells.Item(90, 12).Value = 299997  # L90: 0 -> 299997
$ws.Cells.Item(90, 14).Value = -312477  # N90: None -> -312477
$ws.Cells.Item(107, 8).Value = 1964.625  # H107: 2135.4092 -> 1964.625
$ws.Cells.Item(107, 9).Value = 979.3684  # I107: 1084.4706 -> 979.3684
$ws.Cells.Item(107, 11).Value = 2938.1052  # K107: 3253.4118 -> 2938.1052
$ws.Cells.Item(107, 13).Value = -1018.1052  # M107: -1333.4118 -> -1018.1052
$ws.Cells.Item(122, 8).Value = 2791.5715  # H122: 2858.2341 -> 2791.5715
$ws.Cells.Item(122, 9).Value = 2454.9512  # I122: 2518.0256 -> 2454.9512
$ws.Cells.Item(122, 11).Value = 7364.8536  # K122: 7554.0768 -> 7364.8536
$ws.Cells.Item(122, 13).Value = -4914.8536  # M122: -5104.0768 -> -4914.8536
$ws.Cells.Item(132, 8).Value = 2926250.5  # H132: 3032553 -> 2926250.5
$ws.Cells.Item(132, 9).Value = 3402736.8  # I132: 3402743 -> 3402736.8
$ws.Cells.Item(132, 10).Value = 7772.375  # J132: 9335.666999999999 -> 7772.375
$ws.Cells.Item(132, 11).Value = 10208210.4  # K132: 10208229 -> 10208210.4
$ws.Cells.Item(132, 12).Value = 23317.125  # L132: 28007.001 -> 23317.125
$ws.Cells.Item(132, 13).Value = -10205680.4  # M132: -10205699 -> -10205680.4
$ws.Cells.Item(132, 14).Value = -28377.125  # N132: -33067.001 -> -28377.125
$ws.Cells.Item(136, 8).Value = 23085.715  # H136: 1251741.1 -> 23085.715
$ws.Cells.Item(136, 9).Value = 6097.4  # I136: 1507170.1 -> 6097.4
$ws.Cells.Item(136, 10).Value = 65556.5  # J136: 76768.2 -> 65556.5
$ws.Cells.Item(136, 11).Value = 18292.2  # K136: 4521510.300000001 -> 18292.2
$ws.Cells.Item(136, 12).Value = 196669.5  # L136: 230304.6 -> 196669.5
$ws.Cells.Item(136, 13).Value = -15742.2  # M136: -4518960.300000001 -> -15742.2
$ws.Cells.Item(136, 14).Value = -201769.5  # N136: -235404.6 -> -201769.5
